$d = $word.ActiveDocument

# The source edit drops three whole paragraphs that used to sit right
# after the "LOB1012: Estatistica (Requisito fraco)" paragraph:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and
#       Github pages. Original theme under Creative Commons Attribution"
#
# Locate the anchor paragraph first.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "LOB1012: Estatística (Requisito fraco)*") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    # Only proceed if the three paragraphs being targeted for removal
    # are actually the ones we expect (empty / "Ver no Jupiter..." /
    # "... Creative Commons Attribution"), so this is a no-op if the
    # document was already edited.
    $p1 = $anchor.Next()
    $p2 = $null
    $p3 = $null
    if ($p1 -ne $null) { $p2 = $p1.Next() }
    if ($p2 -ne $null) { $p3 = $p2.Next() }

    $matches = ($p1 -ne $null) -and ($p2 -ne $null) -and ($p3 -ne $null) `
        -and ($p1.Range.Text.Trim() -eq "") `
        -and ($p2.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") `
        -and ($p3.Range.Text -like "*Creative Commons Attribution*")

    if ($matches) {
        # Walk forward three paragraphs from the anchor, collecting the
        # paragraph-mark-inclusive end of the last one, so the deleted
        # range lines up exactly on paragraph boundaries (no leftover /
        # merged paragraph marks).
        $startPos = $anchor.Range.End
        $endPos = $p3.Range.End

        if ($endPos -gt $startPos) {
            $r = $d.Range($startPos, $endPos)
            $r.Delete()
        }
    }
}
